# 5.2.1.1b worksheet: add a new "2023" column (R) mirroring the existing
# year columns, copying formatting from column Q (the previous last year
# column, 2022) and then filling in the 2023 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q formatting (styles) into the new column R for every row
# that currently has data (rows 3-25); this brings across borders/number
# formats/etc. exactly as Excel's own "insert column to the right, copy
# format" flow would.
for ($r = 3; $r -le 25; $r++) {
    $srcCell = $ws.Cells.Item($r, 17)   # column Q
    $dstCell = $ws.Cells.Item($r, 18)   # column R
    $srcCell.Copy($dstCell)
}

# New 2023 data for the numeric rows.
$ws.Range("R4").Value = 2023
$ws.Range("R5").Value = 11357
$ws.Range("R7").Value = 11002
$ws.Range("R8").Value = 355
